# Rename the first sheet from VERIFY_TABLE_DATA to Food_Item
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Food_Item"

# Replace the old USERNAME/PASSWORD/T1/P1/T2/P2 table (columns A & B) with a
# single "Food Item" column (A2:A4) and drop column B's contents entirely.
$ws.Range("A2").Value = "Food Item"
$ws.Range("A3").Value = "Fried Rice"
$ws.Range("A4").Value = "Chicken Tandoori"
$ws.Range("B2:B4").ClearContents()

# Update the active selection to match the new authored state.
$null = $ws.Range("C28").Select()
